{"js": "// Update the date line and all 25 \"three-digit x one-digit\" multiplication\n// answers in the practice-sheet table to the new values from the latest\n// commit (output generated at c8c62b6).\n\nconst replacements = [\n    [\"2025-12-08 Monday\", \"2025-12-09 Tuesday\"],\n    [\"263\u00d75=1315\", \"889\u00d78=7112\"],\n    [\"461\u00d72=922\", \"747\u00d74=2988\"],\n    [\"767\u00d75=3835\", \"776\u00d73=2328\"],\n    [\"751\u00d75=3755\", \"720\u00d78=5760\"],\n    [\"794\u00d74=3176\", \"302\u00d73=906\"],\n    [\"858\u00d76=5148\", \"454\u00d76=2724\"],\n    [\"487\u00d79=4383\", \"965\u00d75=4825\"],\n    [\"562\u00d76=3372\", \"184\u00d76=1104\"],\n    [\"804\u00d78=6432\", \"819\u00d76=4914\"],\n    [\"415\u00d74=1660\", \"895\u00d79=8055\"],\n    [\"304\u00d77=2128\", \"754\u00d73=2262\"],\n    [\"836\u00d75=4180\", \"589\u00d78=4712\"],\n    [\"460\u00d74=1840\", \"200\u00d76=1200\"],\n    [\"857\u00d79=7713\", \"321\u00d72=642\"],\n    [\"651\u00d79=5859\", \"361\u00d75=1805\"],\n    [\"240\u00d74=960\", \"280\u00d76=1680\"],\n    [\"649\u00d79=5841\", \"490\u00d73=1470\"],\n    [\"390\u00d74=1560\", \"677\u00d74=2708\"],\n    [\"112\u00d74=448\", \"749\u00d73=2247\"],\n    [\"460\u00d72=920\", \"424\u00d72=848\"],\n    [\"538\u00d76=3228\", \"306\u00d76=1836\"],\n    [\"831\u00d79=7479\", \"143\u00d72=286\"],\n    [\"176\u00d73=528\", \"594\u00d74=2376\"],\n    [\"476\u00d79=4284\", \"175\u00d79=1575\"],\n    [\"751\u00d73=2253\", \"305\u00d74=1220\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, \"Replace\");\n    }\n    await context.sync();\n}\n", "ps1": "# Update the date line and all 25 \"three-digit x one-digit\" multiplication\n# answers in the practice-sheet table to the new values from the latest\n# commit (output generated at c8c62b6).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-12-08 Monday\"; New = \"2025-12-09 Tuesday\" },\n    @{ Old = \"263\u00d75=1315\";        New = \"889\u00d78=7112\" },\n    @{ Old = \"461\u00d72=922\";         New = \"747\u00d74=2988\" },\n    @{ Old = \"767\u00d75=3835\";        New = \"776\u00d73=2328\" },\n    @{ Old = \"751\u00d75=3755\";        New = \"720\u00d78=5760\" },\n    @{ Old = \"794\u00d74=3176\";        New = \"302\u00d73=906\" },\n    @{ Old = \"858\u00d76=5148\";        New = \"454\u00d76=2724\" },\n    @{ Old = \"487\u00d79=4383\";        New = \"965\u00d75=4825\" },\n    @{ Old = \"562\u00d76=3372\";        New = \"184\u00d76=1104\" },\n    @{ Old = \"804\u00d78=6432\";        New = \"819\u00d76=4914\" },\n    @{ Old = \"415\u00d74=1660\";        New = \"895\u00d79=8055\" },\n    @{ Old = \"304\u00d77=2128\";        New = \"754\u00d73=2262\" },\n    @{ Old = \"836\u00d75=4180\";        New = \"589\u00d78=4712\" },\n    @{ Old = \"460\u00d74=1840\";        New = \"200\u00d76=1200\" },\n    @{ Old = \"857\u00d79=7713\";        New = \"321\u00d72=642\" },\n    @{ Old = \"651\u00d79=5859\";        New = \"361\u00d75=1805\" },\n    @{ Old = \"240\u00d74=960\";         New = \"280\u00d76=1680\" },\n    @{ Old = \"649\u00d79=5841\";        New = \"490\u00d73=1470\" },\n    @{ Old = \"390\u00d74=1560\";        New = \"677\u00d74=2708\" },\n    @{ Old = \"112\u00d74=448\";         New = \"749\u00d73=2247\" },\n    @{ Old = \"460\u00d72=920\";         New = \"424\u00d72=848\" },\n    @{ Old = \"538\u00d76=3228\";        New = \"306\u00d76=1836\" },\n    @{ Old = \"831\u00d79=7479\";        New = \"143\u00d72=286\" },\n    @{ Old = \"176\u00d73=528\";         New = \"594\u00d74=2376\" },\n    @{ Old = \"476\u00d79=4284\";        New = \"175\u00d79=1575\" },\n    @{ Old = \"751\u00d73=2253\";        New = \"305\u00d74=1220\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
